$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the ingredient lists (shared strings content update) for existing products.
$ws.Range("C2").Value = "1.0-Vainilla,2.0-Harinita,1.0-Huevos,5.0-Leche,"
$ws.Range("C3").Value = "2.0-Harinita,5.0-Huevos,"
$ws.Range("C4").Value = "1.0-Vainilla,5.0-Harinita,2.0-Huevos,"
$ws.Range("C5").Value = "5.0-Harinita,5.0-Huevos,"
$ws.Range("C6").Value = "5.0-Crema,4.0-Harinita,5.0-Huevos,2.0-Limon,"

# Add a new order/product row (historial de pedidos de un cliente).
$ws.Range("A8").Value = "test"
$ws.Range("B8").Value = 2.0
$ws.Range("C8").Value = "1.0-Crema,"
$ws.Range("D8").Value = 1000.0
